# RotJ - 596 frame improvement
# Updates the "FrameCounts" sheet (Andymac / Break-last-rock / Joker sequence
# frame counts) to reflect the improved run's new splits, removing rows that
# no longer apply and adding the new Joker-ground split.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 ("1st Hit"): start-of-split time no longer recorded -> remove B74
# entirely (Clear, not ClearContents, so the cell element itself is dropped).
$ws.Range("B74").Clear()

# Row 75 ("Joker on ground (HP = 22)"): start-of-split time no longer recorded -> remove B75.
$ws.Range("B75").Clear()

# Row 76 ("END"): new, faster start time for this split.
$ws.Range("B76").Value = 13269

# Rows 78-80: these intermediate splits (1st Hit / HP=50 / HP=42 on the old
# route) are gone entirely on the improved route - remove place name and both
# time columns.
$ws.Range("A78:C78").Clear()
$ws.Range("A79:C79").Clear()
$ws.Range("A80:C80").Clear()

# Row 82 ("Joker ground 1"): place name and end time no longer used.
$ws.Range("A82").Clear()
$ws.Range("C82").Clear()

# Row 89: new split - Joker ground, start time only recorded so far.
$ws.Range("B89").Value = 12563

# Update the active selection to match the new point of interest (B75).
$ws.Range("B75").Select()
